# Insert a new weekly price record for "Repollo" (Macroferia Regional de Talca)
# just before the old row 430, shifting the existing rows 430-446 down to 431-447.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 430..446 down by one to make room for the new record.
$ws.Rows.Item(430).Insert()

# Populate the newly inserted row 430 with the new weekly observation.
$ws.Cells.Item(430, 1).Value  = 5
$ws.Cells.Item(430, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(430, 3).Value  = "Maule"
$ws.Cells.Item(430, 4).Value  = 45041
$ws.Cells.Item(430, 5).Value  = 7
$ws.Cells.Item(430, 6).Value  = 100112006
$ws.Cells.Item(430, 7).Value  = "Repollo"
$ws.Cells.Item(430, 8).Value  = "Crespo record"
$ws.Cells.Item(430, 9).Value  = "Primera"
$ws.Cells.Item(430, 10).Value = 2000
$ws.Cells.Item(430, 11).Value = 1000
$ws.Cells.Item(430, 12).Value = 1000
$ws.Cells.Item(430, 13).Value = 1000
$ws.Cells.Item(430, 14).Value = "$/unidad"
$ws.Cells.Item(430, 15).Value = "Región del Maule"
$ws.Cells.Item(430, 16).Value = 1000
$ws.Cells.Item(430, 17).Value = 1
$ws.Cells.Item(430, 18).Value = "Hortaliza"
